$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix fuzzed / corrupted raw values
$ws.Range("D14").Value = -457552811.06
$ws.Range("D15").Value = 30000000

# Restore formulas that were flattened to plain values by the fuzzer
$ws.Range("D18").Formula = "=SUM(D12:D17)"
$ws.Range("D21").Formula = "=SUM(D18:D20)"

# Force full recalculation so dependent formula cells (D23, D25, D28, E28, C29, ...)
# pick up the corrected values
$excel.CalculateFullRebuild()
